$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet (values updated first to preserve shared-string order) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokeM"
$wsCuenta.Range("B2").Value = "SmokeLastM"
$wsCuenta.Range("C2").Value = 20100100
$wsCuenta.Range("D2").Value = 103

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 623

# --- DatosMotor sheet ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA004"
$wsMotor.Range("B2").Value = "ABC12SSMA004"
$wsMotor.Range("C2").Value = "ZAZ123SSMA004"

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200103
$wsAP.Activate()
$wsAP.Range("D6").Select()

# --- Re-activate DatosCuenta last so it remains the selected tab/selection ---
$wsCuenta.Activate()
$wsCuenta.Range("D3").Select()
